$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old row 2 was a sub-header row ("Hiver"/"Eté"/"Année"/"(MW)"/"(GWh)" ...).
# Remove it entirely; data below shifts up by one row.
$ws.Rows("2").Delete()

# Rewrite row 1 with the new column headers.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
# E1 held the leftover "Eté" style (s=1) from the deleted sub-header row;
# clear it back to the default/no style like its A1:D1 neighbours.
$ws.Range("E1").Style = "Normal"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 get a dedicated header style (Arial 9, general number format) that is
# distinct from the plain cells in A1:E1. Build it as a transient named style
# so the resulting cell format (applyFont only, no applyNumberFormat) matches
# what a non-Excel writer would emit, then drop the named style itself so the
# workbook's cellStyles/cellStyleXfs tables stay at their original size.
$tmpStyle = $wb.Styles.Add("__tmp_header_style__")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "__tmp_header_style__"
$wb.Styles.Item("__tmp_header_style__").Delete()

# Match the author's resulting selection.
$ws.Range("A2:K2").Select() | Out-Null
